$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the formatting of the
# existing header cell H1 so the new cells share the same cell style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I0 / IF numeric data for rows 2-67
$data = @(
    @{Row=2; I=6; J=6},
    @{Row=3; I=8; J=8},
    @{Row=4; I=9; J=9},
    @{Row=5; I=9; J=9},
    @{Row=6; I=7; J=7},
    @{Row=7; I=5; J=6},
    @{Row=8; I=9; J=9},
    @{Row=9; I=6; J=6},
    @{Row=10; I=8; J=8},
    @{Row=11; I=7; J=8},
    @{Row=12; I=6; J=6},
    @{Row=13; I=5; J=6},
    @{Row=14; I=9; J=10},
    @{Row=15; I=5; J=6},
    @{Row=16; I=7; J=7},
    @{Row=17; I=8; J=8},
    @{Row=18; I=7; J=7},
    @{Row=19; I=8; J=8},
    @{Row=20; I=8; J=8},
    @{Row=21; I=7; J=7},
    @{Row=22; I=8; J=8},
    @{Row=23; I=8; J=8},
    @{Row=24; I=6; J=6},
    @{Row=25; I=8; J=8},
    @{Row=26; I=8; J=8},
    @{Row=27; I=9; J=9},
    @{Row=28; I=6; J=6},
    @{Row=29; I=4; J=5},
    @{Row=30; I=6; J=6},
    @{Row=31; I=7; J=7},
    @{Row=32; I=8; J=8},
    @{Row=33; I=6; J=6},
    @{Row=34; I=1; J=1},
    @{Row=35; I=7; J=7},
    @{Row=36; I=9; J=9},
    @{Row=37; I=5; J=6},
    @{Row=38; I=6; J=6},
    @{Row=39; I=11; J=11},
    @{Row=40; I=8; J=8},
    @{Row=41; I=6; J=6},
    @{Row=42; I=5; J=6},
    @{Row=43; I=4; J=4},
    @{Row=44; I=7; J=7},
    @{Row=45; I=7; J=7},
    @{Row=46; I=7; J=7},
    @{Row=47; I=8; J=8},
    @{Row=48; I=7; J=7},
    @{Row=49; I=9; J=9},
    @{Row=50; I=7; J=8},
    @{Row=51; I=6; J=6},
    @{Row=52; I=7; J=7},
    @{Row=53; I=6; J=7},
    @{Row=54; I=6; J=7},
    @{Row=55; I=7; J=7},
    @{Row=56; I=5; J=5},
    @{Row=57; I=7; J=7},
    @{Row=58; I=6; J=7},
    @{Row=59; I=8; J=8},
    @{Row=60; I=6; J=6},
    @{Row=61; I=6; J=6},
    @{Row=62; I=3; J=4},
    @{Row=63; I=6; J=6},
    @{Row=64; I=5; J=5},
    @{Row=65; I=8; J=8},
    @{Row=66; I=6; J=6},
    @{Row=67; I=3; J=3}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
